{"js": "// Remove the trailing \"Ver no Jupiter...\" line, the \"\u00a9 2020 ...\" copyright\n// line, and the now-superfluous blank paragraph that separated them from the\n// rest of the document, while leaving the single blank paragraph right after\n// \"LOM3099: Est\u00e1tica (Requisito)\" (and the final page-break paragraph) intact.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Ver no Jupiter ...\" paragraph - this is where the block of\n// paragraphs to remove starts.\nlet startIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Ver no Jupiter\") !== -1) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex !== -1) {\n  // Delete: the blank paragraph immediately preceding \"Ver no Jupiter...\"\n  // (startIndex - 1), the \"Ver no Jupiter...\" paragraph itself (startIndex),\n  // and the \"\u00a9 2020 ...\" paragraph right after it (startIndex + 1).\n  const toDelete = [];\n  if (startIndex - 1 >= 0 && items[startIndex - 1].text === \"\") {\n    toDelete.push(items[startIndex - 1]);\n  }\n  toDelete.push(items[startIndex]);\n  if (startIndex + 1 < items.length && items[startIndex + 1].text.indexOf(\"\u00a9\") !== -1) {\n    toDelete.push(items[startIndex + 1]);\n  }\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" line, the \"\u00a9 2020 ...\" copyright\n# line, and the blank paragraph that separated them from the rest of the\n# document, while leaving the blank paragraph right after\n# \"LOM3099: Est\u00e1tica (Requisito)\" (and the final page-break paragraph) intact.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Ver no Jupiter\")\n\nif ($found) {\n    $para = $rng.Paragraphs(1)\n\n    $prevPara = $para.Previous()\n    $nextPara = $para.Next()\n\n    $startPos = $para.Range.Start\n    if ($prevPara -ne $null -and $prevPara.Range.Text.Trim() -eq \"\") {\n        $startPos = $prevPara.Range.Start\n    }\n\n    $endPos = $para.Range.End\n    if ($nextPara -ne $null -and $nextPara.Range.Text -like \"*Contact: luizeleno@usp.br*\") {\n        $endPos = $nextPara.Range.End\n    }\n\n    $delRange = $d.Range($startPos, $endPos)\n    $delRange.Delete()\n}\n"}
